# Apply the sampledata.xlsx edits described by the commit:
#  "Removed taxes page, updated Ratio handling, capital gains tax is now
#   funcitonal and can be found on portfolio page"
#
# Data changes on Sheet1:
#  - Row 5  (MSFT trade on 43494): Type BUY -> SELL, Amount 33.5321 -> 2.1852
#  - Row 18 (MSFT trade on 43507): Amount 8.2939000000000007 -> 1.8239000000000001
#  - Row 28 (FB trade on 43517):   Amount 3.7 -> 0.7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: Type column (C) BUY -> SELL, Amount column (D) 33.5321 -> 2.1852
$ws.Range("C5").Value = "SELL"
$ws.Range("D5").Value = 2.1852

# Row 18: Amount column (D) 8.2939000000000007 -> 1.8239000000000001
$ws.Range("D18").Value = 1.8239000000000001

# Row 28: Amount column (D) 3.7 -> 0.7
$ws.Range("D28").Value = 0.7

# Match the updated active-cell selection recorded in the saved view state.
$ws.Range("J8").Select()
